# harvest_summary.xlsx update:
#  - A new harvest entry ("Swiss Chard", index 25) is inserted before the
#    existing "Onion" row, pushing the Onion/Garlic rows (and the blank
#    spacer + Total row) down by one row.
#  - Several quantity (column E) values are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (shifts Onion, Garlic, blank row, and Total row down by one)
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the Swiss Chard entry
$ws.Cells.Item(9, 1).Value = 25
$ws.Cells.Item(9, 2).Value = "Swiss Chard"
$ws.Cells.Item(9, 3).Value = "Ruby Red"
$ws.Cells.Item(9, 4).Value = "Bunches"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 7).Formula = "=E9*F9"

# Updated quantities (column E) for several existing entries
$ws.Range("E2").Value = 80.62
$ws.Range("E3").Value = 47.52
$ws.Range("E5").Value = 33.61
$ws.Range("E6").Value = 2
$ws.Range("E8").Value = 3.19
